# feat(unit test): add expected mean standard errors to (failing) unit test
# This script adds a new "EXPECTED MEANS" crosstab block (rows 42-48) to
# the worksheet, mirroring the existing EXPECTED PERCENTS / EXPECTED COUNTS
# blocks, to support the crosstab_mean() function unit test.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section title ---
$ws.Range("A42").Value = "EXPECTED MEANS"

# --- Header row (43) ---
$ws.Range("A43").Value = "HHINCOME_bucket"
$ws.Range("B43").Value = "AGE_bucket"
$ws.Range("C43").Value = "RACE_ETH_bucket"
$ws.Range("D43").Value = "weighted_count"
$ws.Range("E43").Value = "wc1"
$ws.Range("F43").Value = "wc2"
$ws.Range("G43").Value = "wc3"
$ws.Range("H43").Value = "wc4"
$ws.Range("I43").Value = "count"
$ws.Range("J43").Value = "weighted_mean"
$ws.Range("K43").Value = "wm1"
$ws.Range("L43").Value = "wm2"
$ws.Range("M43").Value = "wm3"
$ws.Range("N43").Value = "wm4"
$ws.Range("O43").Value = "sqdiff1"
$ws.Range("P43").Value = "sqdiff2"
$ws.Range("Q43").Value = "sqdiff3"
$ws.Range("R43").Value = "sqdiff4"
$ws.Range("S43").Value = "se"

# --- Row 44: r000_100k / r00_49 / white ---
$ws.Range("A44").Value = "r000_100k"
$ws.Range("B44").Value = "r00_49"
$ws.Range("C44").Value = "white"
$ws.Range("D44").Value = 65
$ws.Range("E44").Value = 64
$ws.Range("F44").Value = 71
$ws.Range("G44").Value = 60
$ws.Range("H44").Value = 64
$ws.Range("I44").Value = 2
$ws.Range("J44").Formula = '=(F2*$E2+F14*$E14)/D44'
$ws.Range("K44").Formula = '=(G2*$E2+G14*$E14)/E44'
$ws.Range("L44").Formula = '=(H2*$E2+H14*$E14)/F44'
$ws.Range("M44").Formula = '=(I2*$E2+I14*$E14)/G44'
$ws.Range("N44").Formula = '=(J2*$E2+J14*$E14)/H44'
$ws.Range("O44").Formula = '=(K44-$J44)^2'
$ws.Range("P44").Formula = '=(L44-$J44)^2'
$ws.Range("Q44").Formula = '=(M44-$J44)^2'
$ws.Range("R44").Formula = '=(N44-$J44)^2'
$ws.Range("S44").Formula = '=SQRT((4/80)*SUM(O44:R44))'

# --- Row 45: r000_100k / r00_49 / black ---
$ws.Range("A45").Value = "r000_100k"
$ws.Range("B45").Value = "r00_49"
$ws.Range("C45").Value = "black"
$ws.Range("D45").Value = 116
$ws.Range("E45").Value = 84
$ws.Range("F45").Value = 110
$ws.Range("G45").Value = 111
$ws.Range("H45").Value = 109
$ws.Range("I45").Value = 2
$ws.Range("J45").Formula = '=(F4*$E4+F12*$E12)/D45'
$ws.Range("K45").Formula = '=(G4*$E4+G12*$E12)/E45'
$ws.Range("L45").Formula = '=(H4*$E4+H12*$E12)/F45'
$ws.Range("M45").Formula = '=(I4*$E4+I12*$E12)/G45'
$ws.Range("N45").Formula = '=(J4*$E4+J12*$E12)/H45'
$ws.Range("O45").Formula = '=(K45-$J45)^2'
$ws.Range("P45").Formula = '=(L45-$J45)^2'
$ws.Range("Q45").Formula = '=(M45-$J45)^2'
$ws.Range("R45").Formula = '=(N45-$J45)^2'
$ws.Range("S45").Formula = '=SQRT((4/80)*SUM(O45:R45))'

# --- Row 46: r000_100k / r50plus / black ---
$ws.Range("A46").Value = "r000_100k"
$ws.Range("B46").Value = "r50plus"
$ws.Range("C46").Value = "black"
$ws.Range("D46").Value = 106
$ws.Range("E46").Value = 129
$ws.Range("F46").Value = 106
$ws.Range("G46").Value = 103
$ws.Range("H46").Value = 102
$ws.Range("I46").Value = 3
$ws.Range("J46").Formula = '=((F3*$E3)+($E10*F10)+(F15*$E15))/D46'
$ws.Range("K46").Formula = '=((G3*$E3)+($E10*G10)+(G15*$E15))/E46'
$ws.Range("L46").Formula = '=((H3*$E3)+($E10*H10)+(H15*$E15))/F46'
$ws.Range("M46").Formula = '=((I3*$E3)+($E10*I10)+(I15*$E15))/G46'
$ws.Range("N46").Formula = '=((J3*$E3)+($E10*J10)+(J15*$E15))/H46'
$ws.Range("O46").Formula = '=(K46-$J46)^2'
$ws.Range("P46").Formula = '=(L46-$J46)^2'
$ws.Range("Q46").Formula = '=(M46-$J46)^2'
$ws.Range("R46").Formula = '=(N46-$J46)^2'
$ws.Range("S46").Formula = '=SQRT((4/80)*SUM(O46:R46))'

# --- Row 47: r000_100k / r50plus / aian ---
$ws.Range("A47").Value = "r000_100k"
$ws.Range("B47").Value = "r50plus"
$ws.Range("C47").Value = "aian"
$ws.Range("D47").Value = 99
$ws.Range("E47").Value = 95
$ws.Range("F47").Value = 96
$ws.Range("G47").Value = 107
$ws.Range("H47").Value = 98
$ws.Range("I47").Value = 2
$ws.Range("J47").Formula = '=(F11*$E11+F13*$E13)/D47'
$ws.Range("K47").Formula = '=(G11*$E11+G13*$E13)/E47'
$ws.Range("L47").Formula = '=(H11*$E11+H13*$E13)/F47'
$ws.Range("M47").Formula = '=(I11*$E11+I13*$E13)/G47'
$ws.Range("N47").Formula = '=(J11*$E11+J13*$E13)/H47'
$ws.Range("O47").Formula = '=(K47-$J47)^2'
$ws.Range("P47").Formula = '=(L47-$J47)^2'
$ws.Range("Q47").Formula = '=(M47-$J47)^2'
$ws.Range("R47").Formula = '=(N47-$J47)^2'
$ws.Range("S47").Formula = '=SQRT((4/80)*SUM(O47:R47))'

# --- Row 48: r100kplus / r00_49 / aapi ---
$ws.Range("A48").Value = "r100kplus"
$ws.Range("B48").Value = "r00_49"
$ws.Range("C48").Value = "aapi"
$ws.Range("D48").Value = 228
$ws.Range("E48").Value = 217
$ws.Range("F48").Value = 227
$ws.Range("G48").Value = 255
$ws.Range("H48").Value = 222
$ws.Range("I48").Value = 5
$ws.Range("J48").Formula = '=SUMPRODUCT(F5:F9,$E5:$E9)/D48'
$ws.Range("K48").Formula = '=SUMPRODUCT(G5:G9,$E5:$E9)/E48'
$ws.Range("L48").Formula = '=SUMPRODUCT(H5:H9,$E5:$E9)/F48'
$ws.Range("M48").Formula = '=SUMPRODUCT(I5:I9,$E5:$E9)/G48'
$ws.Range("N48").Formula = '=SUMPRODUCT(J5:J9,$E5:$E9)/H48'
$ws.Range("O48").Formula = '=(K48-$J48)^2'
$ws.Range("P48").Formula = '=(L48-$J48)^2'
$ws.Range("Q48").Formula = '=(M48-$J48)^2'
$ws.Range("R48").Formula = '=(N48-$J48)^2'
$ws.Range("S48").Formula = '=SQRT((4/80)*SUM(O48:R48))'

# --- Column widths for the new columns (best-effort; engine quantizes to 1/6) ---
$ws.Columns("D").ColumnWidth = 16.45
$ws.Columns("E").ColumnWidth = 8.88
$ws.Columns("F").ColumnWidth = 15.02
$ws.Columns("G").ColumnWidth = 19.02
$ws.Columns("J").ColumnWidth = 13.74

# --- View state: scroll/selection to match final cursor position ---
$ws.Range("S47").Select()
